$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 205, shifting the existing rows 205:216 down to 206:217.
$ws.Rows("205:205").Insert()

# Fill in the new row 205 with the new weekly record.
$ws.Range("A205").Value = 10
$ws.Range("B205").Value = "Vega Modelo de Temuco"
$ws.Range("C205").Value = "La Araucanía"
$ws.Range("D205").Value = 44578
$ws.Range("E205").Value = 9
$ws.Range("F205").Value = "Fruta"
$ws.Range("G205").Value = 100103
$ws.Range("H205").Value = "Frutos de hueso (carozo)"
$ws.Range("I205").Value = 100103001
$ws.Range("J205").Value = "Cereza"
$ws.Range("K205").Value = "Lapins"
$ws.Range("L205").Value = "Primera"
$ws.Range("M205").Value = 650
$ws.Range("N205").Value = 8000
$ws.Range("O205").Value = 8000
$ws.Range("P205").Value = 8000
$ws.Range("Q205").Value = "$/bandeja 10 kilos"
$ws.Range("R205").Value = "Región de La Araucanía"
$ws.Range("S205").Value = 800
$ws.Range("T205").Value = 10
